$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")
$ws.Select()

$ws.Range("H1").Value = "Sort"
$ws.Range("H2").Value = "lohi"
$ws.Range("H3").Value = "hilo"

$ws.Range("H1:H3").HorizontalAlignment = -4108

$ws.Range("H3").Select()
